$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for the new "label" column
$ws.Range("B1").Value = "label"

# Fill B2:B11 with "humano" label value
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = "humano"
}
